# Update the cryptocurrency price/volume list with freshly scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new text value. All of these are plain-text cells in the
# original workbook (inline strings), including ones that look numeric
# (e.g. "22.18") or percentages (e.g. "  -1.46%  "). Writing a number-looking
# string straight into .Value lets Excel auto-coerce it into a real number,
# so for every cell we briefly force Text number-format ("@") before the
# assignment, then reset the style back to "Normal" so no stray formatting
# is left behind (matching the original, unstyled inline-string cells).
function Set-TextValue {
    param($range, [string]$value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$textUpdates = [ordered]@{
    "D2"  = "27.164.51"
    "E2"  = "  -1.46%  "
    "D3"  = "1.572.51"
    "E3"  = "  -0.80%  "
    "E4"  = "  +0.02%  "
    "D5"  = "207.60"
    "E5"  = "  -0.09%  "
    "D6"  = "0.490"
    "E6"  = "  -1.96%  "
    "E7"  = "  +0.07%  "
    "D8"  = "22.18"
    "E8"  = "  -0.27%  "
    "D9"  = "0.248"
    "E9"  = "  -0.93%  "
    "D10" = "0.0590"
    "E10" = "  +0.12%  "
    "E11" = "  -0.02%  "
    "D12" = "1.800.71"
    "E12" = "  -0.57%  "
    "D13" = "1.575.62"
    "E13" = "  -0.30%  "
    "D14" = "3.77"
    "E14" = "  -1.42%  "
    "D15" = "0.518"
    "E15" = "  -1.18%  "
    "D16" = "27.215.00"
    "E16" = "  -1.33%  "
    "D17" = "62.17"
    "E17" = "  -1.32%  "
    "D18" = "214.21"
    "E18" = "  -1.01%  "
    "E19" = "  +0.62%  "
    "D20" = "0.0₃0686"
    "E20" = "  -0.87%  "
    "E21" = "  -0.06%  "
    "D22" = "4.12"
    "E22" = "  -0.31%  "
    "D23" = "9.41"
    "E23" = "  -3.03%  "
    "E24" = "  +0.53%  "
    "D25" = "151.77"
    "E25" = "  -0.67%  "
    "D26" = "6.65"
    "E26" = "  -4.10%  "
    "D27" = "14.95"
    "E27" = "  -0.52%  "
    "E28" = "  +0.02%  "
    "E29" = "  -1.13%  "
    "E30" = "  -1.65%  "
    "D31" = "0.0463"
    "E31" = "  -1.94%  "
    "D32" = "3.17"
    "E32" = "  -1.43%  "
    "D33" = "1.406.41"
    "E33" = "  +2.46%  "
    "E34" = "  -1.25%  "
    "D35" = "1.55"
    "E35" = "  +1.37%  "
    "E36" = "  -0.86%  "
    "D37" = "0.938"
    "E37" = "  -2.97%  "
    "D38" = "0.0165"
    "E38" = "  -1.70%  "
    "D39" = "0.817"
    "E39" = "  -0.67%  "
    "D40" = "0.516"
    "E40" = "  -3.05%  "
    "E41" = "  +0.06%  "
    "D42" = "0.994"
    "E42" = "  +2.31%  "
    "E43" = "  +3.45%  "
    "E44" = "  +2.17%  "
    "D45" = "2.18"
    "E45" = "  +0.48%  "
    "D46" = "63.66"
    "E46" = "  -0.76%  "
    "D47" = "1.710.49"
    "E47" = "  -0.73%  "
    "D48" = "86.02"
    "E48" = "  -0.45%  "
    "D49" = "0.0₇0981"
    "E49" = "  -2.60%  "
}

foreach ($addr in $textUpdates.Keys) {
    Set-TextValue $ws.Range($addr) $textUpdates[$addr]
}

# Rows 50/51: Algorand and Cronos swap places, with refreshed price/volume figures.
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D50") "0.0494"
$ws.Range("E50").Value = "  +0.07%  "

$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws.Range("D51") "0.0951"
$ws.Range("E51").Value = "  -0.77%  "
